$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the measured results for MQTT QoS1 (row 4) and MQTT QoS2 (row 5) ---
$ws.Range("B4").Value = 3.0176198963727199
$ws.Range("C4").Value = 0.88701854723590701
$ws.Range("D4").Value = 178.10329826233399
$ws.Range("E4").Value = 33.119611438055003
$ws.Range("F4").Value = 1624.72417478068
$ws.Range("G4").Value = 820.59410309504005
$ws.Range("H4").Value = 811.87211008102702
$ws.Range("I4").Value = 167.52533726585
$ws.Range("J4").Value = 1.20999999999974
$ws.Range("K4").Value = 1.02620000000001
$ws.Range("L4").Value = 1.0240224609374999
$ws.Range("M4").Value = 1.0078306640624899

$ws.Range("B5").Value = 1.03959373775608
$ws.Range("C5").Value = 0.293787289958928
$ws.Range("D5").Value = 79.943300124757499
$ws.Range("E5").Value = 15.188980261646799
$ws.Range("F5").Value = 1517.32661261556
$ws.Range("G5").Value = 456.79227781619699
$ws.Range("H5").Value = 1928.6923804854
$ws.Range("I5").Value = 506.44985189132899
$ws.Range("J5").Value = 1.20999999999974
$ws.Range("K5").Value = 1.02620000000001
$ws.Range("L5").Value = 1.0240224609374999
$ws.Range("M5").Value = 1.0078306640624899

# --- Apply the 3-decimal numeric display format to the whole results block ---
$ws.Range("B4:M7").NumberFormat = "0.000"

# --- Update the active selection left on the sheet ---
$ws.Range("H2:I2").Select() | Out-Null
